# Replace every (non-overlapping) occurrence of $Find inside a TextRange's
# text with $ReplaceWith, working right-to-left through the match positions
# so earlier replacements never invalidate the character offsets used by
# later ones. Uses TextRange.Characters(start, length) sub-ranges so the
# run-level formatting (rPr) of whichever run each match sits in is kept.
function Replace-InTextRange($TextRange, $Find, $ReplaceWith) {
    $full = $TextRange.Text
    $positions = @()
    $idx = $full.IndexOf($Find)
    while ($idx -ge 0) {
        $positions += $idx
        $idx = $full.IndexOf($Find, $idx + 1)
    }
    for ($i = $positions.Count - 1; $i -ge 0; $i--) {
        $pos = $positions[$i] + 1
        $sub = $TextRange.Characters($pos, $Find.Length)
        $sub.Text = $ReplaceWith
    }
    return $positions.Count
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 ("Remove"): min = none / max = none -> null
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(3)
$tr10 = $shp10.TextFrame.TextRange
Replace-InTextRange $tr10 "min = none" "min = null" | Out-Null
Replace-InTextRange $tr10 "    max = none" "    max = null" | Out-Null

# ---------------------------------------------------------------------
# Slide 11 ("Successor"): USize -> universeSize, none -> null (incl.
# "!=" / "==" comparisons)
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(3)
$tr11 = $shp11.TextFrame.TextRange
Replace-InTextRange $tr11 "USize" "universeSize" | Out-Null
Replace-InTextRange $tr11 "min != none " "min != null " | Out-Null
Replace-InTextRange $tr11 "max != none " "max != null " | Out-Null
Replace-InTextRange $tr11 "temp == none" "temp == null" | Out-Null
Replace-InTextRange $tr11 "none" "null" | Out-Null

# ---------------------------------------------------------------------
# Slide 12 ("Predecessor"): same kind of changes as slide 11
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(3)
$tr12 = $shp12.TextFrame.TextRange
Replace-InTextRange $tr12 "USize" "universeSize" | Out-Null
Replace-InTextRange $tr12 "min != none " "min != null " | Out-Null
Replace-InTextRange $tr12 "max != none " "max != null " | Out-Null
Replace-InTextRange $tr12 "temp == none" "temp == null" | Out-Null
Replace-InTextRange $tr12 "none" "null" | Out-Null

# ---------------------------------------------------------------------
# Slide 9 ("Insert"): move/resize the body placeholder, lower-case the
# "If " runs, and apply the same USize/none -> universeSize/null swaps.
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(3)

# Width must be set before Left (setting Left first freezes Width at its
# old value in this host), and 453.42855pt is the value that round-trips
# to the exact target EMU extent (5758542).
$shp9.Width = 453.42855
$shp9.Left = 66.0

$tr9 = $shp9.TextFrame.TextRange
Replace-InTextRange $tr9 "    If " "    if " | Out-Null
Replace-InTextRange $tr9 "If " "if " | Out-Null
Replace-InTextRange $tr9 "USize" "universeSize" | Out-Null
Replace-InTextRange $tr9 "none" "null" | Out-Null
